# "Price Lookup" sheet: pick a Product/Brand in the lookup cells and
# highlight where that price lives in the table above, via conditional
# formatting keyed off the looked-up result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Choose a product/brand combination in the lookup cells (B12/C12) so the
# INDEX/MATCH formula in D12 resolves to a real price instead of the
# IFERROR fallback of 0.
$ws.Range("B12").Value = "Mouse"
$ws.Range("C12").Value = "Lenovo"

# Re-enter D12's formula so it is re-evaluated against the new B12/C12
# inputs and its cached result is refreshed.
$ws.Range("D12").Formula = "=IFERROR(INDEX(B2:F6,MATCH(B12,A2:A6,0),MATCH(C12,B1:F1,0)),0)"

# B2 picks up the same currency/fill formatting already used by the rest
# of row 2 (C2:F2) instead of its own separate (but visually identical)
# cell format.
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Conditionally highlight every cell in the price table (B2:F7) whose
# "Apple" column value equals the price looked up in D12, so the matching
# row/column stands out. Interior.Color is an OLE (BGR) value; 49407
# (0x00C0FF, bytes B=00 G=C0 R=FF) is RGB #FFC000 - a golden/amber fill.
$rng = $ws.Range("B2:F7")
$fc = $rng.FormatConditions.Add(2, 3, 'B2=$D$12')
$fc.Interior.Color = 49407

# Reflect the final selection landing on C12 (the brand dropdown) rather
# than D12.
[void]$ws.Range("C12").Select()
